$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '96.652.51'
$ws.Range("E2").Value = '  +1.96%  '

$ws.Range("D3").Value = '3.278.10'
$ws.Range("E3").Value = '  +4.77%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '249.00'
$ws.Range("E5").Value = '  +4.40%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '617.00'
$ws.Range("E6").Value = '  +0.24%  '

$ws.Range("E7").Value = '  -1.33%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.380'
$ws.Range("E8").Value = '  -3.22%  '

$ws.Range("E9").Value = '  +0.03%  '

$ws.Range("D10").Value = '3.268.69'
$ws.Range("E10").Value = '  +4.58%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.778'
$ws.Range("E11").Value = '  -6.84%  '

$ws.Range("E12").Value = '  -0.02%  '

$ws.Range("D13").Value = '96.313.43'
$ws.Range("E13").Value = '  +2.00%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000243'
$ws.Range("E14").Value = '  -1.04%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '34.97'
$ws.Range("E15").Value = '  +0.32%  '

$ws.Range("D16").Value = '3.874.95'
$ws.Range("E16").Value = '  +4.71%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.47'
$ws.Range("E17").Value = '  +3.04%  '

$ws.Range("D18").Value = '3.269.78'
$ws.Range("E18").Value = '  +4.60%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.55'
$ws.Range("E19").Value = '  -4.13%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.81'
$ws.Range("E20").Value = '  -1.25%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '476.49'
$ws.Range("E21").Value = '  +5.76%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.76'
$ws.Range("E22").Value = '  -2.73%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.0000203'
$ws.Range("E23").Value = '  +1.66%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.14'
$ws.Range("E24").Value = '  +1.23%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.57'
$ws.Range("E25").Value = '  -1.31%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '87.39'
$ws.Range("E26").Value = '  +1.55%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.95'
$ws.Range("E27").Value = '  -1.77%  '

$ws.Range("D28").Value = '3.453.72'
$ws.Range("E28").Value = '  +4.90%  '

$ws.Range("E29").Value = '  +0.16%  '

$ws.Range("E30").Value = '  -1.44%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.237'
$ws.Range("E31").Value = '  -9.74%  '

$ws.Range("E32").Value = '  +0.99%  '

$ws.Range("E33").Value = '  -2.21%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '9.14'
$ws.Range("E34").Value = '  -2.93%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '27.09'
$ws.Range("E35").Value = '  +3.97%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '7.32'
$ws.Range("E36").Value = '  -8.57%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.149'
$ws.Range("E37").Value = '  -6.82%  '

$ws.Range("E38").Value = '  +0.24%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '24.70'
$ws.Range("E39").Value = '  +3.00%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '490.46'
$ws.Range("E40").Value = '  +2.38%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.442'
$ws.Range("E41").Value = '  -3.52%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.26'
$ws.Range("E42").Value = '  -3.83%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.22'
$ws.Range("E43").Value = '  -0.35%  '

$ws.Range("B44").Value = 'ARBITRUM'
$ws.Range("C44").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.783'
$ws.Range("E44").Value = '  +12.59%  '

$ws.Range("B45").Value = 'USDe'
$ws.Range("C45").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.00'
$ws.Range("E45").Value = '  +0.00%  '

$ws.Range("B46").Value = 'MantraDAO'
$ws.Range("C46").Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.43'
$ws.Range("E46").Value = '  -5.43%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '161.06'
$ws.Range("E47").Value = '  -0.38%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.88'
$ws.Range("E48").Value = '  -0.45%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '45.13'
$ws.Range("E49").Value = '  +2.64%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.34'
$ws.Range("E50").Value = '  +2.12%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.46'
$ws.Range("E51").Value = '  +0.73%  '
